$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell L4: "Multivalued" (bold, Calibri 11 - same family/colour as
# the sheet's base font, just bold + smaller, matching the other headers).
$ws.Range("L4").Value = "Multivalued"
$hf = $ws.Range("L4").Font
$hf.Bold = $true
$hf.Size = 11

# Data cells L5:L7: boolean-looking text "FALSE", displayed through a
# custom TRUE/FALSE number format, left aligned. The leading apostrophe
# forces a real text cell instead of a native boolean, matching the
# existing TRUE/FALSE text columns elsewhere on this sheet.
$ws.Range("L5:L7").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("L5:L7").HorizontalAlignment = -4131
$ws.Range("L5:L7").Value = "'FALSE"

# Rows 4-7 shrink slightly (15pt) to match the new content.
$ws.Range("A4:L7").RowHeight = 15

# Selection moves onto the newly added column.
$ws.Range("L4:L7").Select()
